$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new rows (49, 50, 51) of webcam data ---
# Copy formatting from row 48 (last existing row) onto the new row 49, matching
# the style pattern Excel keeps when a new row is typed right under the table.
$ws.Range("A48:F48").Copy()
$ws.Range("A49:F49").PasteSpecial(-4122)

# Rows 50 and 51 pick up the plainer formatting pattern (no special F-column
# border, default B-column style) seen elsewhere in the table, e.g. row 44.
$ws.Range("A44:F44").Copy()
$ws.Range("A50:F50").PasteSpecial(-4122)
$ws.Range("A44:F44").Copy()
$ws.Range("A51:F51").PasteSpecial(-4122)

# Fill in the values. The order below reproduces the exact order the three
# rows were actually populated in (not a simple left-to-right fill), which
# keeps the shared-string table's append order correct.
$ws.Range("A49").Value = "LIVE, CITY, TRAFFIC, BRIDGE"
$ws.Range("B49").Value = "59.938937786015906, 30.30665733597569"
$ws.Range("C49").Value = "LIVE CAMERA Sidewalk 24/7 St. Petersburg"
$ws.Range("D49").Value = "St. Petersburg"
$ws.Range("E49").Value = "Russia"
$ws.Range("F49").Value = "DGHhKZlFpXM"

$ws.Range("F50").Value = "wULEcVypV4Q"
$ws.Range("D50").Value = "Vladivostok"
$ws.Range("B50").Value = "43.10994206708416, 131.88786535463603"
$ws.Range("C50").Value = "Online Camera Vladivostok Center Веб-камера Владивосток Океанский пр-т"

$ws.Range("F51").Value = "utI0WcV36Tk"
$ws.Range("E51").Value = "Spain"
$ws.Range("B51").Value = "37.97718790981845, -0.6702841528155156"
$ws.Range("D51").Value = "Torrevieja"
$ws.Range("C51").Value = "Playa del Cura"

$ws.Range("A50").Value = "LIVE, CITY, SEA, SHIP"
$ws.Range("A51").Value = "LIVE, SEA, BEACH"
$ws.Range("E50").Value = "Russia"

# --- Apply a fresh "Highlight duplicate values" conditional format to the
# new F49 cell (same red-fill rule already used for the rest of column F). ---
$fc = $ws.Range("F49").FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Interior.Color = 255
$fc.SetFirstPriority()

# --- Reflect the final cursor/scroll position like Excel would after typing
# the last new row. ---
$ws.Range("A51").Select()
$excel.ActiveWindow.ScrollRow = 28
